# Updated cryptos list on Sat Jan  6 22:50:35 UTC 2024 with GitHub Actions
# Applies per-cell updates (price / volume / coin swap) to the crypto tracker sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Cell = "D2"; Value = "44.280.72"; ForceText = $false },
    @{ Cell = "E2"; Value = "  +0.54%  "; ForceText = $false },
    @{ Cell = "D3"; Value = "2.237.35"; ForceText = $false },
    @{ Cell = "E3"; Value = "  -0.16%  "; ForceText = $false },
    @{ Cell = "E4"; Value = "  +0.10%  "; ForceText = $false },
    @{ Cell = "D5"; Value = "306.88"; ForceText = $true },
    @{ Cell = "E5"; Value = "  -2.77%  "; ForceText = $false },
    @{ Cell = "D6"; Value = "94.11"; ForceText = $true },
    @{ Cell = "E6"; Value = "  -4.61%  "; ForceText = $false },
    @{ Cell = "E7"; Value = "  -0.82%  "; ForceText = $false },
    @{ Cell = "D8"; Value = "1.01"; ForceText = $true },
    @{ Cell = "E8"; Value = "  +0.31%  "; ForceText = $false },
    @{ Cell = "E9"; Value = "  -2.16%  "; ForceText = $false },
    @{ Cell = "E10"; Value = "  -4.89%  "; ForceText = $false },
    @{ Cell = "E11"; Value = "  -1.94%  "; ForceText = $false },
    @{ Cell = "D12"; Value = "7.14"; ForceText = $true },
    @{ Cell = "E12"; Value = "  -2.89%  "; ForceText = $false },
    @{ Cell = "E13"; Value = "  -0.11%  "; ForceText = $false },
    @{ Cell = "D14"; Value = "2.578.14"; ForceText = $false },
    @{ Cell = "E14"; Value = "  -0.17%  "; ForceText = $false },
    @{ Cell = "D15"; Value = "2.310.51"; ForceText = $false },
    @{ Cell = "E15"; Value = "  +2.87%  "; ForceText = $false },
    @{ Cell = "D16"; Value = "0.825"; ForceText = $true },
    @{ Cell = "E16"; Value = "  -2.15%  "; ForceText = $false },
    @{ Cell = "D17"; Value = "13.42"; ForceText = $true },
    @{ Cell = "E17"; Value = "  -3.68%  "; ForceText = $false },
    @{ Cell = "D18"; Value = "43.984.57"; ForceText = $false },
    @{ Cell = "E18"; Value = "  +0.22%  "; ForceText = $false },
    @{ Cell = "D19"; Value = "0.0₃0962"; ForceText = $false },
    @{ Cell = "E19"; Value = "  -2.10%  "; ForceText = $false },
    @{ Cell = "E20"; Value = "  +0.35%  "; ForceText = $false },
    @{ Cell = "E21"; Value = "  -8.49%  "; ForceText = $false },
    @{ Cell = "D22"; Value = "65.51"; ForceText = $true },
    @{ Cell = "E22"; Value = "  +0.29%  "; ForceText = $false },
    @{ Cell = "E23"; Value = "  +4.23%  "; ForceText = $false },
    @{ Cell = "D24"; Value = "236.74"; ForceText = $true },
    @{ Cell = "E24"; Value = "  -0.80%  "; ForceText = $false },
    @{ Cell = "E25"; Value = "  -1.49%  "; ForceText = $false },
    @{ Cell = "E26"; Value = "  +0.35%  "; ForceText = $false },
    @{ Cell = "D27"; Value = "39.60"; ForceText = $true },
    @{ Cell = "E27"; Value = "  +6.24%  "; ForceText = $false },
    @{ Cell = "E28"; Value = "  +4.13%  "; ForceText = $false },
    @{ Cell = "D29"; Value = "9.85"; ForceText = $true },
    @{ Cell = "D30"; Value = "20.01"; ForceText = $true },
    @{ Cell = "D31"; Value = "5.84"; ForceText = $true },
    @{ Cell = "E31"; Value = "  -2.63%  "; ForceText = $false },
    @{ Cell = "D32"; Value = "153.19"; ForceText = $true },
    @{ Cell = "E32"; Value = "  -1.56%  "; ForceText = $false },
    @{ Cell = "D33"; Value = "0.0793"; ForceText = $true },
    @{ Cell = "E33"; Value = "  -5.20%  "; ForceText = $false },
    @{ Cell = "E34"; Value = "  -2.58%  "; ForceText = $false },
    @{ Cell = "D35"; Value = "3.08"; ForceText = $true },
    @{ Cell = "E35"; Value = "  -9.93%  "; ForceText = $false },
    @{ Cell = "E36"; Value = "  +1.51%  "; ForceText = $false },
    @{ Cell = "D37"; Value = "0.107"; ForceText = $true },
    @{ Cell = "E37"; Value = "  -3.50%  "; ForceText = $false },
    @{ Cell = "E38"; Value = "  -7.98%  "; ForceText = $false },
    @{ Cell = "D39"; Value = "3.47"; ForceText = $true },
    @{ Cell = "E39"; Value = "  -1.30%  "; ForceText = $false },
    @{ Cell = "D40"; Value = "3.79"; ForceText = $true },
    @{ Cell = "E40"; Value = "  -2.72%  "; ForceText = $false },
    @{ Cell = "D41"; Value = "14.11"; ForceText = $true },
    @{ Cell = "E41"; Value = "  -7.42%  "; ForceText = $false },
    @{ Cell = "E42"; Value = "  -3.77%  "; ForceText = $false },
    @{ Cell = "E43"; Value = "  +0.24%  "; ForceText = $false },
    @{ Cell = "D44"; Value = "1.718.72"; ForceText = $false },
    @{ Cell = "E44"; Value = "  +0.52%  "; ForceText = $false },
    @{ Cell = "D45"; Value = "82.28"; ForceText = $true },
    @{ Cell = "E45"; Value = "  -1.29%  "; ForceText = $false },
    @{ Cell = "E46"; Value = "  -2.11%  "; ForceText = $false },
    @{ Cell = "D47"; Value = "4.92"; ForceText = $true },
    @{ Cell = "E47"; Value = "  -5.05%  "; ForceText = $false },
    @{ Cell = "D48"; Value = "99.06"; ForceText = $true },
    @{ Cell = "E48"; Value = "  -2.43%  "; ForceText = $false },
    @{ Cell = "E49"; Value = "  -0.45%  "; ForceText = $false },
    @{ Cell = "B50"; Value = "MultiversX"; ForceText = $false },
    @{ Cell = "C50"; Value = "https://coinranking.com/coin/omwkOTglq+multiversx-egld"; ForceText = $false },
    @{ Cell = "D50"; Value = "54.58"; ForceText = $true },
    @{ Cell = "E50"; Value = "  -3.28%  "; ForceText = $false },
    @{ Cell = "B51"; Value = "FraxShare"; ForceText = $false },
    @{ Cell = "C51"; Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"; ForceText = $false },
    @{ Cell = "D51"; Value = "8.03"; ForceText = $true },
    @{ Cell = "E51"; Value = "  -0.97%  "; ForceText = $false }
)

foreach ($item in $updates) {
    $rng = $ws.Range($item.Cell)
    if ($item.ForceText) {
        # Force the cell to remain text so Excel doesn't reinterpret
        # decimal-looking values (e.g. "1.01") as numbers.
        $rng.NumberFormat = "@"
        $rng.Value = $item.Value
        $rng.Style = "Normal"
    } else {
        $rng.Value = $item.Value
    }
}
